# "more bulk upload fixes"
# Insert a new "Date Created (Year)*" column after the filename column on
# Sheet1, fill it with 2000 for every existing data row, give that new
# column's data cells an explicit black font color, and leave the
# selection sitting on the new header cell (B1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B (shifts former B:F -> C:G).
$ws.Columns.Item(2).Insert() | Out-Null

# New header + year values for the inserted column.
$ws.Range("B1").Value = "Date Created (Year)*"
$ws.Range("B2:B4").Value = 2000

# Explicit black font color on the newly entered year cells.
$ws.Range("B2:B4").Font.Color = 0

# Leave the active selection on the new header cell.
$ws.Range("B1").Select() | Out-Null
